# Updated cryptos list on Mon Mar 27 11:36:37 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values scraped
# from coinranking.com for each coin row in the sheet.
#
# Column D holds price strings that look numeric (e.g. "329.26"); a leading
# apostrophe forces Excel to store them as text (matching the workbook's
# original inline-string / General-format cells) instead of silently
# converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.979.86"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "'1.769.93"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'329.26"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.4669"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").Value = "'0.3525"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'43.93"
$ws.Range("E9").Value = "  +5.14%  "
$ws.Range("D10").Value = "'0.07394"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "'1.087"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D13").Value = "'20.66"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'6.023"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'7.203"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "'1.767.18"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'92.27"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'0.06418"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'16.95"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "'5.805"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'28.013.12"
$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "'2.156"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").Value = "'164.50"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "'1.969.27"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'2.203"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'123.45"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'0.09336"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").Value = "'3.658"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'5.559"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'11.69"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "'0.02270"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'0.06116"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "'0.2074"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'4.917"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "'1.195"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "'0.6168"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").Value = "'1.447"
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").Value = "'7.786"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "'13.13"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "'3.749"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'0.5812"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'124.08"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'1.938"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'0.06813"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("E51").Value = "  -0.36%  "
